$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table with the latest scraped values.
# All of these cells hold plain text in the workbook (coin names, links,
# price strings and percentage-change strings). Some of the new price
# strings look like ordinary numbers (e.g. "543.84", "1.00", "0.0000160"),
# which Excel would otherwise auto-convert to a numeric value and thereby
# mangle (dropping trailing zeros, switching to scientific notation, etc.).
# A leading apostrophe forces Excel to keep those as literal text, just
# like the original inline-string cells.

$ws.Range("D2").Value = "60.777.71"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.349.38"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'543.84"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'136.34"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Value = "2.348.43"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'24.64"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "2.773.27"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "60.759.44"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'0.0000160"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "2.352.40"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "'10.63"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'319.76"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'4.12"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'6.55"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'63.38"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "'1.66"
$ws.Range("E25").Value = "  -9.47%  "
$ws.Range("D26").Value = "'8.44"
$ws.Range("E26").Value = "  +8.47%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'7.98"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'492.50"
$ws.Range("E30").Value = "  -5.33%  "
$ws.Range("D31").Value = "0.0₃0863"
$ws.Range("E31").Value = "  -7.46%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.146"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.62"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.376"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'18.50"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "  +5.35%  "
$ws.Range("D40").Value = "'5.22"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'145.12"
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'142.08"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.57"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.04"
$ws.Range("E45").Value = "  -7.88%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0515"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'19.04"
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.568"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0900"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0221"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.39"
$ws.Range("E51").Value = "  +0.17%  "
